$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as the first data row of the
# "Apio" (Macroferia Regional de Talca) block. Insert a new row at 113
# (pushing the existing rows 113-130 down to 114-131) and populate it
# with the new observation.
$ws.Rows.Item(113).Insert()

$ws.Range("A113").Value = 5
$ws.Range("B113").Value = "Macroferia Regional de Talca"
$ws.Range("C113").Value = "Maule"
$ws.Range("D113").Value = 44491
$ws.Range("E113").Value = 7
$ws.Range("F113").Value = 100112017
$ws.Range("G113").Value = "Apio"
$ws.Range("H113").Value = "Americana (o)"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 7000
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = 7000
$ws.Range("N113").Value = '$/docena de matas'
$ws.Range("O113").Value = "Provincia del Elquí"
$ws.Range("P113").Value = 1167
$ws.Range("Q113").Value = 6
$ws.Range("R113").Value = "Hortaliza"
